$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8502261638641357
$ws.Range("B1").Value = 2.663928031921387
$ws.Range("C1").Value = 3.230233669281006
$ws.Range("D1").Value = 1.817433834075928
$ws.Range("E1").Value = 1.391499042510986
